# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes to "_FV2210" / "_FV2304"
# - Turn the data range into a real Excel Table ("Table1")
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row labels: columns A-J keep the "_old" base name but get
#    the "_FV2210" suffix, columns L-U (the "_new" variants) get "_FV2304".
#    Column K ("diff") is left untouched.
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $leftCol = $i + 1        # columns 1..10  => A..J
    $rightCol = $i + 12      # columns 12..21 => L..U
    $ws.Cells.Item(1, $leftCol).Value = $baseNames[$i] + "_FV2210"
    $ws.Cells.Item(1, $rightCol).Value = $baseNames[$i] + "_FV2304"
}

# 2) Convert A1:U84 into an Excel table named "Table1"
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U84"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (pane split under row 1)
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
